$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("B3:B10").Formula = "=B1+B2"

[void]$ws.Range("B3:B10").Select()
